$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows (old rows 17 and 18); the table now ends at row 16
$ws.Rows("17:18").Delete() | Out-Null

# Update measured values for rows 2-16 with the refreshed COMPASS pidis deuteron data
$ws.Cells.Item(2,3).Value = 0.0045999999999999999
$ws.Cells.Item(2,4).Value = 1.1000000000000001
$ws.Cells.Item(2,7).Value = -0.0054000000000000003
$ws.Cells.Item(2,8).Value = 0.0074000000000000003
$ws.Cells.Item(2,9).Value = 0.0047999999999999996
$ws.Cells.Item(3,3).Value = 0.0054999999999999997
$ws.Cells.Item(3,4).Value = 1.22
$ws.Cells.Item(3,7).Value = 0.00029999999999999997
$ws.Cells.Item(3,8).Value = 0.0057999999999999996
$ws.Cells.Item(3,9).Value = 0.0043
$ws.Cells.Item(4,3).Value = 0.0070000000000000001
$ws.Cells.Item(4,4).Value = 1.39
$ws.Cells.Item(4,7).Value = -0.0011000000000000001
$ws.Cells.Item(4,8).Value = 0.0041999999999999997
$ws.Cells.Item(4,9).Value = 0.0023
$ws.Cells.Item(5,3).Value = 0.0089999999999999993
$ws.Cells.Item(5,4).Value = 1.62
$ws.Cells.Item(5,7).Value = -0.0086999999999999994
$ws.Cells.Item(5,8).Value = 0.0048999999999999998
$ws.Cells.Item(5,9).Value = 0.0030999999999999999
$ws.Cells.Item(6,3).Value = 0.0141
$ws.Cells.Item(6,4).Value = 2.19
$ws.Cells.Item(6,7).Value = -0.0011000000000000001
$ws.Cells.Item(6,8).Value = 0.0032000000000000002
$ws.Cells.Item(6,9).Value = 0.0023999999999999998
$ws.Cells.Item(7,3).Value = 0.024400000000000002
$ws.Cells.Item(7,4).Value = 3.29
$ws.Cells.Item(7,7).Value = 0.0074999999999999997
$ws.Cells.Item(7,8).Value = 0.0047999999999999996
$ws.Cells.Item(7,9).Value = 0.0033999999999999998
$ws.Cells.Item(8,3).Value = 0.034599999999999999
$ws.Cells.Item(8,4).Value = 4.43
$ws.Cells.Item(8,7).Value = 0.0094999999999999998
$ws.Cells.Item(8,8).Value = 0.0064000000000000003
$ws.Cells.Item(8,9).Value = 0.0041999999999999997
$ws.Cells.Item(9,3).Value = 0.0487
$ws.Cells.Item(9,4).Value = 6.06
$ws.Cells.Item(9,7).Value = 0.015900000000000001
$ws.Cells.Item(9,8).Value = 0.0063
$ws.Cells.Item(9,9).Value = 0.0044000000000000003
$ws.Cells.Item(10,3).Value = 0.076600000000000001
$ws.Cells.Item(10,4).Value = 9
$ws.Cells.Item(10,7).Value = 0.052699999999999997
$ws.Cells.Item(10,8).Value = 0.0070000000000000001
$ws.Cells.Item(10,9).Value = 0.0071999999999999998
$ws.Cells.Item(11,3).Value = 0.121
$ws.Cells.Item(11,4).Value = 13.5
$ws.Cells.Item(11,7).Value = 0.095000000000000001
$ws.Cells.Item(11,8).Value = 0.01
$ws.Cells.Item(11,9).Value = 0.010999999999999999
$ws.Cells.Item(12,3).Value = 0.17100000000000001
$ws.Cells.Item(12,4).Value = 18.600000000000001
$ws.Cells.Item(12,7).Value = 0.121
$ws.Cells.Item(12,8).Value = 0.014999999999999999
$ws.Cells.Item(12,9).Value = 0.016
$ws.Cells.Item(13,3).Value = 0.222
$ws.Cells.Item(13,4).Value = 23.8
$ws.Cells.Item(13,7).Value = 0.16
$ws.Cells.Item(13,8).Value = 0.021999999999999999
$ws.Cells.Item(13,9).Value = 0.02
$ws.Cells.Item(14,3).Value = 0.28999999999999998
$ws.Cells.Item(14,4).Value = 31.1
$ws.Cells.Item(14,7).Value = 0.19
$ws.Cells.Item(14,8).Value = 0.023
$ws.Cells.Item(14,9).Value = 0.021999999999999999
$ws.Cells.Item(15,3).Value = 0.40500000000000003
$ws.Cells.Item(15,4).Value = 43.9
$ws.Cells.Item(15,7).Value = 0.317
$ws.Cells.Item(15,8).Value = 0.036999999999999998
$ws.Cells.Item(15,9).Value = 0.035999999999999997
$ws.Cells.Item(16,3).Value = 0.56699999999999995
$ws.Cells.Item(16,4).Value = 60.8
$ws.Cells.Item(16,7).Value = 0.49399999999999999
$ws.Cells.Item(16,8).Value = 0.082000000000000003
$ws.Cells.Item(16,9).Value = 0.084000000000000005

# Restore the shared formulas: J2/K2 stand alone, J3:J16/K3:K16 form shared groups
$ws.Cells.Item(2,10).Formula = "=G2/10"
$ws.Cells.Item(2,11).Formula = "=MAX(I2^2-J2^2,0)^0.5"
$ws.Range("J3:J16").Formula = "=G3/10"
$ws.Range("K3:K16").Formula = "=MAX(I3^2-J3^2,0)^0.5"

# Update the active selection to match the saved view
$ws.Range("I17").Select() | Out-Null
